$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "Neurolucida results"
$ws.Range("B23").Value = "2022-06-13 15-18-49"
$ws.Range("C23").Value = "DEG"
$ws.Range("D23").Value = "SCTv2 corrected BL_A + BL_C new post selection"
$ws.Range("F23").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G23").Value = "pseudotime"

$ws.Range("G24").Select()
